$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.930.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.553.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("D5").Value = "'207.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.488"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").Value = "'21.76"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.58%  "
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D13").Value = "'1.554.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").Value = "'0.518"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "'26.929.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'61.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "'216.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.04%  "
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "'153.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("D26").Value = "'6.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "'14.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").Value = "'0.0468"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.37%  "
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'1.434.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.37%  "
$ws.Range("D34").Value = "'3.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.67%  "
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").Value = "'2.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "'0.519"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").Value = "'5.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").Value = "'0.991"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("D45").Value = "'63.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("D47").Value = "'1.689.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").Value = "'86.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").Value = "'0.0525"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.23%  "
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("D51").Value = "'0.0957"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.74%  "
